# hw2/Bond.xlsx -- bondPrice function work with simple maturity timeframe
#
# The second cash-flow schedule (rows 16-22) modeled every coupon date the
# same way, including the final one. In reality the last cash flow at
# maturity also repays the principal, so it should just be the hard
# number (coupon + par) rather than the generic coupon formula. Wire that
# up, extend the running "price so far" total through that last column,
# and drop the stray scratch formula that was left below the table.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# J17 = final cash flow at maturity = last coupon (4) + principal redemption (100)
$ws.Range("J17").Value = 104

# Fill out the running total row so it also covers the discounted final
# cash flow (previously stopped one column short of J22).
$ws.Range("I22").Formula = "=SUM(C20:I20)"

# Remove the leftover one-off scratch calculation below the table.
$ws.Range("A25:L25").EntireRow.Delete()

# Leave the selection where the author left off reviewing the change.
$ws.Range("F10").Select()
